$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "55×49=2695" "25×91=2275"
Replace-Text "29×21=609" "66×60=3960"
Replace-Text "46×16=736" "73×32=2336"
Replace-Text "95×15=1425" "45×70=3150"
Replace-Text "86×37=3182" "12×32=384"
Replace-Text "57×38=2166" "64×52=3328"
Replace-Text "26×68=1768" "56×34=1904"
Replace-Text "98×64=6272" "55×58=3190"
Replace-Text "84×94=7896" "80×15=1200"
Replace-Text "33×76=2508" "41×69=2829"
Replace-Text "53×35=1855" "41×64=2624"
Replace-Text "35×14=490" "60×49=2940"
Replace-Text "30×71=2130" "86×88=7568"
Replace-Text "74×60=4440" "38×40=1520"
Replace-Text "73×96=7008" "19×25=475"
Replace-Text "40×69=2760" "51×22=1122"
Replace-Text "19×90=1710" "69×21=1449"
Replace-Text "92×72=6624" "55×96=5280"
Replace-Text "19×45=855" "88×35=3080"
Replace-Text "92×45=4140" "74×71=5254"
Replace-Text "34×13=442" "43×62=2666"
Replace-Text "89×35=3115" "87×16=1392"
Replace-Text "25×78=1950" "58×47=2726"
Replace-Text "42×15=630" "29×52=1508"
Replace-Text "50×57=2850" "26×74=1924"
